$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '64.510.46'
$ws.Range("D2").Style = "Normal"

$ws.Range("E2").Value = '  +2.76%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.526.99'
$ws.Range("D3").Style = "Normal"

$ws.Range("E3").Value = '  +2.74%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '580.87'
$ws.Range("D5").Style = "Normal"

$ws.Range("E5").Value = '  +1.31%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '152.81'
$ws.Range("D6").Style = "Normal"

$ws.Range("E6").Value = '  +4.52%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("D7").Style = "Normal"

$ws.Range("E7").Value = '  -0.06%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.531.94'
$ws.Range("D9").Style = "Normal"

$ws.Range("E9").Value = '  +2.96%  '

$ws.Range("E10").Value = '  +0.77%  '

$ws.Range("E11").Value = '  -1.30%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.30'
$ws.Range("D12").Style = "Normal"

$ws.Range("E12").Value = '  +0.03%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.356'
$ws.Range("D13").Style = "Normal"

$ws.Range("E13").Value = '  -0.52%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '29.43'
$ws.Range("D14").Style = "Normal"

$ws.Range("E14").Value = '  +1.58%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000181'
$ws.Range("D15").Style = "Normal"

$ws.Range("E15").Value = '  +2.17%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.974.39'
$ws.Range("D16").Style = "Normal"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '63.993.93'
$ws.Range("D17").Style = "Normal"

$ws.Range("E17").Value = '  +2.04%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.531.60'
$ws.Range("D18").Style = "Normal"

$ws.Range("E18").Value = '  +2.79%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.87'
$ws.Range("D19").Style = "Normal"

$ws.Range("E19").Value = '  -1.29%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.00'
$ws.Range("D20").Style = "Normal"

$ws.Range("E20").Value = '  +0.25%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.27'
$ws.Range("D21").Style = "Normal"

$ws.Range("E21").Value = '  +3.23%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '328.55'
$ws.Range("D22").Style = "Normal"

$ws.Range("E22").Value = '  +0.68%  '

$ws.Range("E23").Value = '  +1.50%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.999'
$ws.Range("D24").Style = "Normal"

$ws.Range("E24").Value = '  +0.12%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '10.04'
$ws.Range("D25").Style = "Normal"

$ws.Range("E25").Value = '  -0.88%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '65.55'
$ws.Range("D26").Style = "Normal"

$ws.Range("E26").Value = '  -0.06%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '652.09'
$ws.Range("D27").Style = "Normal"

$ws.Range("E27").Value = '  -0.63%  '

$ws.Range("E28").Value = '  +5.70%  '

$ws.Range("E30").Value = '  +5.10%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.01'
$ws.Range("D31").Style = "Normal"

$ws.Range("E31").Value = '  +0.65%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.04'
$ws.Range("D32").Style = "Normal"

$ws.Range("E32").Value = '  +0.68%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.88'
$ws.Range("D33").Style = "Normal"

$ws.Range("E33").Value = '  +1.75%  '

$ws.Range("E34").Value = '  +2.09%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.997'
$ws.Range("D35").Style = "Normal"

$ws.Range("E35").Value = '  -0.04%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.56'
$ws.Range("D36").Style = "Normal"

$ws.Range("E36").Value = '  +1.76%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.83'
$ws.Range("D37").Style = "Normal"

$ws.Range("E37").Value = '  +1.76%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.56'
$ws.Range("D38").Style = "Normal"

$ws.Range("E38").Value = '  +3.19%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.373'
$ws.Range("D39").Style = "Normal"

$ws.Range("E39").Value = '  +1.12%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '18.94'
$ws.Range("D40").Style = "Normal"

$ws.Range("E40").Value = '  +1.39%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '152.76'
$ws.Range("D41").Style = "Normal"

$ws.Range("E41").Value = '  +0.74%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.82'
$ws.Range("D42").Style = "Normal"

$ws.Range("E42").Value = '  +1.41%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.79'
$ws.Range("D43").Style = "Normal"

$ws.Range("E43").Value = '  +3.21%  '

$ws.Range("E44").Value = '  +1.03%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '162.42'
$ws.Range("D45").Style = "Normal"

$ws.Range("E45").Value = '  +5.98%  '

$ws.Range("E47").Value = '  -2.19%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '15.46'
$ws.Range("D48").Style = "Normal"

$ws.Range("E48").Value = '  +1.56%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.65'
$ws.Range("D49").Style = "Normal"

$ws.Range("E49").Value = '  +2.10%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '21.36'
$ws.Range("D50").Style = "Normal"

$ws.Range("E50").Value = '  +4.51%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.620'
$ws.Range("D51").Style = "Normal"

$ws.Range("E51").Value = '  +2.25%  '
